$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84 - this shifts the existing rows 84-95 down to 85-96,
# preserving their data and formatting (including the date style on column D).
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Cells.Item(84, 1).Value2 = 3
$ws.Cells.Item(84, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(84, 3).Value2 = "Coquimbo"
$ws.Cells.Item(84, 4).Value2 = 45127
$ws.Cells.Item(84, 5).Value2 = 5
$ws.Cells.Item(84, 6).Value2 = 100112022
$ws.Cells.Item(84, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(84, 8).Value2 = "Perfection"
$ws.Cells.Item(84, 9).Value2 = "Primera"
$ws.Cells.Item(84, 10).Value2 = 35
$ws.Cells.Item(84, 11).Value2 = 28000
$ws.Cells.Item(84, 12).Value2 = 28000
$ws.Cells.Item(84, 13).Value2 = 28000
$ws.Cells.Item(84, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(84, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(84, 16).Value2 = 1120
$ws.Cells.Item(84, 17).Value2 = 25
$ws.Cells.Item(84, 18).Value2 = "Hortaliza"
